$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.038.32"
$ws.Range("E2").Value = "  -1.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.515.59"
$ws.Range("E3").Value = "  -2.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.78"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.613"
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.503.35"
$ws.Range("E8").Value = "  -2.69%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("E11").Value = "  +2.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.582"
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "47.41"
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000275"
$ws.Range("E14").Value = "  -2.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.079.63"
$ws.Range("E15").Value = "  -2.55%  "
$ws.Range("E16").Value = "  -3.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "627.53"
$ws.Range("E17").Value = "  -5.88%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.087.07"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.500.46"
$ws.Range("E19").Value = "  -2.94%  "
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.43"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.16"
$ws.Range("E22").Value = "  -1.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.888"
$ws.Range("E23").Value = "  -4.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.92"
$ws.Range("E24").Value = "  -6.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "97.17"
$ws.Range("E25").Value = "  -2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.82"
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.64"
$ws.Range("E28").Value = "  -4.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.27"
$ws.Range("E29").Value = "  -6.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.74"
$ws.Range("E30").Value = "  -4.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.56"
$ws.Range("E31").Value = "  -4.11%  "
$ws.Range("E32").Value = "  -5.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.33"
$ws.Range("E33").Value = "  -3.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.94"
$ws.Range("E34").Value = "  -5.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "637.27"
$ws.Range("E35").Value = "  +10.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.76"
$ws.Range("E36").Value = "  -2.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.50"
$ws.Range("E37").Value = "  -10.80%  "
$ws.Range("E38").Value = "  -3.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "57.19"
$ws.Range("E39").Value = "  -1.25%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.135"
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.389.55"
$ws.Range("E43").Value = "  -4.92%  "
$ws.Range("E44").Value = "  -4.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "32.82"
$ws.Range("E45").Value = "  -4.77%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0699"
$ws.Range("E46").Value = "  -4.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.55"
$ws.Range("E47").Value = "  -4.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.74"
$ws.Range("E48").Value = "  -4.85%  "
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "132.25"
$ws.Range("E50").Value = "  -2.59%  "
$ws.Range("E51").Value = "  +14.83%  "
